$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column (D) keeps being stored as plain text, exactly
# like the source feed renders it (it already used text cells, e.g.
# "26.495.89", "1.0000", etc). Without forcing the Text number format,
# Excel would auto-convert clean numeric-looking strings to real numbers.
# Restoring the style to "Normal" afterwards keeps the cell's visual
# formatting/style index identical to the original (unstyled) cell while
# the stored cell content remains text.
$priceCells = "D2","D3","D5","D7","D8","D9","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D50"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.489.19"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "1.851.47"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "261.45"
$ws.Range("E5").Value = "  -7.44%  "

$ws.Range("D7").Value = "0.5147"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "0.3274"
$ws.Range("E8").Value = "  -7.32%  "

$ws.Range("D9").Value = "0.06779"
$ws.Range("E9").Value = "  -4.56%  "

$ws.Range("E10").Value = "  -5.48%  "

$ws.Range("D11").Value = "0.7751"
$ws.Range("E11").Value = "  -5.62%  "

$ws.Range("D12").Value = "0.07700"
$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("D13").Value = "1.868.38"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").Value = "88.70"
$ws.Range("E14").Value = "  -0.96%  "

$ws.Range("D15").Value = "5.044"
$ws.Range("E15").Value = "  -2.11%  "

$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("E17").Value = "  -2.06%  "

$ws.Range("D18").Value = "0.9999"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "0.000007916"
$ws.Range("E19").Value = "  -3.49%  "

$ws.Range("D20").Value = "26.512.62"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").Value = "2.080.25"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").Value = "4.555"
$ws.Range("E22").Value = "  -4.96%  "

$ws.Range("D23").Value = "9.555"
$ws.Range("E23").Value = "  -5.94%  "

$ws.Range("D24").Value = "5.959"
$ws.Range("E24").Value = "  -4.68%  "

$ws.Range("D25").Value = "2.367"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("D26").Value = "144.64"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").Value = "1.658"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("E28").Value = "  -2.70%  "

$ws.Range("D29").Value = "111.28"

$ws.Range("D30").Value = "4.217"
$ws.Range("E30").Value = "  -4.33%  "

$ws.Range("D31").Value = "4.170"
$ws.Range("E31").Value = "  -4.44%  "

$ws.Range("D32").Value = "0.08759"
$ws.Range("E32").Value = "  -0.88%  "

$ws.Range("D33").Value = "0.04859"
$ws.Range("E33").Value = "  -1.12%  "

$ws.Range("E34").Value = "  -3.32%  "

$ws.Range("D35").Value = "2.838"
$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("D36").Value = "0.6925"
$ws.Range("E36").Value = "  -7.47%  "

$ws.Range("D37").Value = "3.122"
$ws.Range("E37").Value = "  -5.14%  "

$ws.Range("D38").Value = "0.01808"
$ws.Range("E38").Value = "  -3.99%  "

$ws.Range("D39").Value = "2.222"
$ws.Range("E39").Value = "  -8.84%  "

$ws.Range("D40").Value = "0.4931"
$ws.Range("E40").Value = "  -6.87%  "

$ws.Range("D41").Value = "113.54"
$ws.Range("E41").Value = "  -2.70%  "

$ws.Range("D42").Value = "0.8974"
$ws.Range("E42").Value = "  -7.87%  "

$ws.Range("D43").Value = "6.148"
$ws.Range("E43").Value = "  -2.51%  "

# Rows 44 and 45 swap places: Aptos <-> PaxDollar
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "7.809"
$ws.Range("E45").Value = "  -4.73%  "

$ws.Range("D46").Value = "0.4268"
$ws.Range("E46").Value = "  -7.50%  "

$ws.Range("D47").Value = "0.1270"
$ws.Range("E47").Value = "  -6.97%  "

$ws.Range("D48").Value = "9.119"
$ws.Range("E48").Value = "  -4.35%  "

$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").Value = "34.98"
$ws.Range("E50").Value = "  -4.42%  "

$ws.Range("E51").Value = "  -6.16%  "

# Restore the original (default/unstyled) cell style on the price cells now
# that the values are committed as text, so formatting matches the source.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
